$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column X (shifts nothing right of it; extends the used range)
# so that formatting is inherited from the adjacent column W, matching style s=7
$ws.Columns("X").Insert()

# Header cell X1: build the literal text via a formula (so Excel does not
# auto-parse the date-like string into a date serial), then convert the
# formula down to a plain value in place, preserving the inherited style.
$ws.Range("X1").Formula = '="09-10-2020"'
$ws.Range("X1").Copy()
$ws.Range("X1").PasteSpecial(-4163)

# Data cells X2:X36: plain numeric values
$ws.Range("X2").Value = 3707
$ws.Range("X3").Value = 684930
$ws.Range("X4").Value = 8679
$ws.Range("X5").Value = 159836
$ws.Range("X6").Value = 180696
$ws.Range("X7").Value = 11344
$ws.Range("X8").Value = 106027
$ws.Range("X9").Value = 3025
$ws.Range("X10").Value = 272948
$ws.Range("X11").Value = 31902
$ws.Range("X12").Value = 127786
$ws.Range("X13").Value = 126267
$ws.Range("X14").Value = 13597
$ws.Range("X15").Value = 69020
$ws.Range("X16").Value = 80439
$ws.Range("X17").Value = 552519
$ws.Range("X18").Value = 167256
$ws.Range("X19").Value = 3540
$ws.Range("X20").Value = 122687
$ws.Range("X21").Value = 1212016
$ws.Range("X22").Value = 9719
$ws.Range("X23").Value = 4832
$ws.Range("X24").Value = 1937
$ws.Range("X25").Value = 5564
$ws.Range("X26").Value = 216984
$ws.Range("X27").Value = 25256
$ws.Range("X28").Value = 107200
$ws.Range("X29").Value = 129618
$ws.Range("X30").Value = 2650
$ws.Range("X31").Value = 586454
$ws.Range("X32").Value = 180953
$ws.Range("X33").Value = 23474
$ws.Range("X34").Value = 44808
$ws.Range("X35").Value = 378662
$ws.Range("X36").Value = 249737

$excel.CutCopyMode = $false
